# Removed parenthetical info from 'Number of People Living with HIV' by using
# find and replace with regexp ']*' (1 operations, 59 records affected).
#
# Effectively: for every data row, strip everything starting at the first
# "[" character (the "[low - high]" range and any trailing parenthetical
# note) from the "Number of People Living with HIV" column (column C),
# leaving just the leading point estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 58 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $orig = $cell.Value2
    if ($orig -ne $null) {
        $text = [string]$orig
        $idx = $text.IndexOf("[")
        if ($idx -ge 0) {
            $cell.Value = $text.Substring(0, $idx)
        }
    }
}

# Restore the active selection recorded in the saved workbook.
$ws.Range("G5").Select() | Out-Null
